# Auto update stock data
# Updates the report date (column A) from 2025/12/20 -> 2025/12/21 for every
# company block, and the EBITDA value (column B) for the first company
# (Alcoa, row 2) from 6.31 -> 6.37, matching the source data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "'2025/12/21"
}

$b2 = $ws.Cells.Item(2, 2)
$b2.Value = "'6.37"
